# Regenerated synthetic dataset sheet: the data-gen script grew multi-index /
# multi-table support, so this sample now spans 20 weekly date columns
# (A:T, up from A:L) with refreshed random values, a couple of header gaps
# (style kept, value blank) and the "soon" placeholder column moved/duplicated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Phase 1: copy style-2 format (bold+border+center, no number format) from K1
#     to the cells that need it in the new layout, before K1 itself is repurposed.
$ws.Range("K1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("S1").PasteSpecial(-4122) | Out-Null
$ws.Range("T1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Phase 2: copy style-1 format (date format, bold+border+center) from A1
#     to the new/repurposed date cells.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null
$ws.Range("O1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").PasteSpecial(-4122) | Out-Null
$ws.Range("Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("R1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Phase 3: clear cells that are removed in the new layout (incl. emptied,
#     still-styled header cells whose old value must not linger)
$ws.Range("B5").ClearContents() | Out-Null
$ws.Range("D3").ClearContents() | Out-Null
$ws.Range("F5").ClearContents() | Out-Null
$ws.Range("G1").ClearContents() | Out-Null
$ws.Range("I4").ClearContents() | Out-Null
$ws.Range("I5").ClearContents() | Out-Null
$ws.Range("S1").ClearContents() | Out-Null

# --- Phase 4: set date values (row 1, style-1 cells)
$ws.Range("A1").Value = 45940
$ws.Range("B1").Value = 45947
$ws.Range("C1").Value = 45954
$ws.Range("D1").Value = 45961
$ws.Range("E1").Value = 45968
$ws.Range("F1").Value = 45975
$ws.Range("H1").Value = 45989
$ws.Range("I1").Value = 45996
$ws.Range("J1").Value = 46003
$ws.Range("K1").Value = 46010
$ws.Range("M1").Value = 46024
$ws.Range("N1").Value = 46031
$ws.Range("O1").Value = 46038
$ws.Range("P1").Value = 46045
$ws.Range("Q1").Value = 46052
$ws.Range("R1").Value = 46059

# --- Phase 5: set the shared "soon" text on the style-2 text cells
$ws.Range("L1").Value = "soon"
$ws.Range("T1").Value = "soon"

# --- Phase 6: set plain numeric data values (rows 2-5, etc.)
$ws.Range("A2").Value = 183
$ws.Range("B2").Value = 127
$ws.Range("C2").Value = 196
$ws.Range("D2").Value = 179
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 67
$ws.Range("G2").Value = 143
$ws.Range("H2").Value = 185
$ws.Range("I2").Value = 126
$ws.Range("J2").Value = 194
$ws.Range("K2").Value = 70
$ws.Range("L2").Value = 124
$ws.Range("M2").Value = 178
$ws.Range("N2").Value = 92
$ws.Range("O2").Value = 50
$ws.Range("P2").Value = 190
$ws.Range("Q2").Value = 154
$ws.Range("R2").Value = 122
$ws.Range("S2").Value = 120
$ws.Range("T2").Value = 180
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 2
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 4
$ws.Range("M3").Value = 9
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 9
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 6
$ws.Range("A4").Value = 142
$ws.Range("B4").Value = 174
$ws.Range("C4").Value = 125
$ws.Range("D4").Value = 188
$ws.Range("E4").Value = 173
$ws.Range("F4").Value = 159
$ws.Range("G4").Value = 137
$ws.Range("H4").Value = 166
$ws.Range("J4").Value = 159
$ws.Range("K4").Value = 58
$ws.Range("L4").Value = 88
$ws.Range("M4").Value = 169
$ws.Range("N4").Value = 56
$ws.Range("O4").Value = 132
$ws.Range("P4").Value = 148
$ws.Range("Q4").Value = 117
$ws.Range("R4").Value = 85
$ws.Range("S4").Value = 66
$ws.Range("T4").Value = 177
$ws.Range("A5").Value = 6
$ws.Range("D5").Value = 3
$ws.Range("G5").Value = 10
$ws.Range("K5").Value = 8
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 8
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 5
$ws.Range("P5").Value = 2
$ws.Range("Q5").Value = 10
$ws.Range("R5").Value = 10
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 7
